# Auto-generated edit script: updates market-price derived columns (H-N)
# on each profession sheet, per the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 336
$ws.Range("I18").Value = 181.66667
$ws.Range("K18").Value = 181.66667
$ws.Range("M18").Value = 102.33333
$ws.Range("H64").Value = 2881.25
$ws.Range("I64").Value = 2877.7778
$ws.Range("K64").Value = 2877.7778
$ws.Range("M64").Value = -2629.7778
$ws.Range("H67").Value = 2881.25
$ws.Range("I67").Value = 2877.7778
$ws.Range("K67").Value = 2877.7778
$ws.Range("M67").Value = -2019.7778
$ws.Range("H69").Value = 13378.75
$ws.Range("I69").Value = 7833.3335
$ws.Range("J69").Value = 30015
$ws.Range("K69").Value = 23500.0005
$ws.Range("L69").Value = 90045
$ws.Range("M69").Value = -22626.0005
$ws.Range("N69").Value = -91793
$ws.Range("H70").Value = 2555.375
$ws.Range("I70").Value = 872.8570999999999
$ws.Range("J70").Value = 3864
$ws.Range("K70").Value = 2618.5713
$ws.Range("L70").Value = 11592
$ws.Range("M70").Value = -2348.5713
$ws.Range("N70").Value = -12132
$ws.Range("H72").Value = 13378.75
$ws.Range("I72").Value = 7833.3335
$ws.Range("J72").Value = 30015
$ws.Range("K72").Value = 70500.0015
$ws.Range("L72").Value = 270135
$ws.Range("M72").Value = -66132.0015
$ws.Range("N72").Value = -278871
$ws.Range("H73").Value = 2555.375
$ws.Range("I73").Value = 872.8570999999999
$ws.Range("J73").Value = 3864
$ws.Range("K73").Value = 2618.5713
$ws.Range("L73").Value = 11592
$ws.Range("M73").Value = -1682.5713
$ws.Range("N73").Value = -13464
$ws.Range("H76").Value = 3422.5
$ws.Range("I76").Value = 3500
$ws.Range("K76").Value = 3500
$ws.Range("M76").Value = -3185
$ws.Range("H79").Value = 3422.5
$ws.Range("I79").Value = 3500
$ws.Range("K79").Value = 3500
$ws.Range("M79").Value = -2408
$ws.Range("H80").Value = 2125
$ws.Range("J80").Value = 2166.6667
$ws.Range("L80").Value = 6500.000100000001
$ws.Range("N80").Value = -8496.000100000001
$ws.Range("H83").Value = 2125
$ws.Range("J83").Value = 2166.6667
$ws.Range("L83").Value = 19500.0003
$ws.Range("N83").Value = -29484.0003
$ws.Range("H137").Value = 2089.7144
$ws.Range("I137").Value = 1229.6666
$ws.Range("J137").Value = 7250
$ws.Range("K137").Value = 3688.9998
$ws.Range("L137").Value = 21750
$ws.Range("M137").Value = -1138.9998
$ws.Range("N137").Value = -26850
$ws.Range("H138").Value = 3468.9
$ws.Range("I138").Value = 2724.8572
$ws.Range("J138").Value = 3567.17
$ws.Range("K138").Value = 8174.571599999999
$ws.Range("L138").Value = 10701.51
$ws.Range("M138").Value = -3034.571599999999
$ws.Range("N138").Value = -20981.51

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1923.6571
$ws.Range("I74").Value = 1510.8928
$ws.Range("J74").Value = 3574.7144
$ws.Range("K74").Value = 1510.8928
$ws.Range("L74").Value = 3574.7144
$ws.Range("M74").Value = -636.8928000000001
$ws.Range("N74").Value = -5322.7144
$ws.Range("H77").Value = 1923.6571
$ws.Range("I77").Value = 1510.8928
$ws.Range("J77").Value = 3574.7144
$ws.Range("K77").Value = 7554.464
$ws.Range("L77").Value = 17873.572
$ws.Range("M77").Value = -3186.464
$ws.Range("N77").Value = -26609.572

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H102").Value = 7846.2856
$ws.Range("I102").Value = 7846.2856
$ws.Range("K102").Value = 7846.2856
$ws.Range("M102").Value = -4601.2856
$ws.Range("H105").Value = 2301.225
$ws.Range("J105").Value = 2147.6428
$ws.Range("L105").Value = 2147.6428
$ws.Range("N105").Value = -5641.6428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2136.4482
$ws.Range("I58").Value = 1757
$ws.Range("J58").Value = 5425
$ws.Range("K58").Value = 1757
$ws.Range("L58").Value = 5425
$ws.Range("M58").Value = -1554
$ws.Range("N58").Value = -5831
$ws.Range("H136").Value = 2136.4482
$ws.Range("I136").Value = 1757
$ws.Range("J136").Value = 5425
$ws.Range("K136").Value = 5271
$ws.Range("L136").Value = 16275
$ws.Range("M136").Value = -2721
$ws.Range("N136").Value = -21375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 400
$ws.Range("I92").Value = 400
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 1200
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 48
$ws.Range("N92").ClearContents()
$ws.Range("H137").Value = 6116.5
$ws.Range("J137").Value = 7988.6665
$ws.Range("L137").Value = 23965.9995
$ws.Range("N137").Value = -34165.99950000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 35716844
$ws.Range("I80").Value = 125001150
$ws.Range("J80").Value = 3120
$ws.Range("K80").Value = 125001150
$ws.Range("L80").Value = 3120
$ws.Range("M80").Value = -125000152
$ws.Range("N80").Value = -5116
$ws.Range("H83").Value = 35716844
$ws.Range("I83").Value = 125001150
$ws.Range("J83").Value = 3120
$ws.Range("K83").Value = 625005750
$ws.Range("L83").Value = 15600
$ws.Range("M83").Value = -625000758
$ws.Range("N83").Value = -25584
$ws.Range("H107").Value = 4831525
$ws.Range("I107").Value = 496.6875
$ws.Range("J107").Value = 15873876
$ws.Range("K107").Value = 496.6875
$ws.Range("L107").Value = 15873876
$ws.Range("M107").Value = 1423.3125
$ws.Range("N107").Value = -15877716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 5062.3794
$ws.Range("I82").Value = 7811.143
$ws.Range("J82").Value = 2496.8667
$ws.Range("K82").Value = 7811.143
$ws.Range("L82").Value = 2496.8667
$ws.Range("M82").Value = -7450.143
$ws.Range("N82").Value = -3218.8667
$ws.Range("H85").Value = 5062.3794
$ws.Range("I85").Value = 7811.143
$ws.Range("J85").Value = 2496.8667
$ws.Range("K85").Value = 7811.143
$ws.Range("L85").Value = 2496.8667
$ws.Range("M85").Value = -6563.143
$ws.Range("N85").Value = -4992.8667
$ws.Range("H93").Value = 6175155.5
$ws.Range("I93").Value = 11112760
$ws.Range("J93").Value = 3149.75
$ws.Range("K93").Value = 11112760
$ws.Range("L93").Value = 3149.75
$ws.Range("M93").Value = -11111512
$ws.Range("N93").Value = -5645.75
$ws.Range("H136").Value = 3133.7585
$ws.Range("I136").Value = 1435.909
$ws.Range("J136").Value = 8469.857
$ws.Range("K136").Value = 4307.727000000001
$ws.Range("L136").Value = 25409.571
$ws.Range("M136").Value = -1757.727000000001
$ws.Range("N136").Value = -30509.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 281001.75
$ws.Range("I62").Value = 12002
$ws.Range("J62").Value = 550001.5
$ws.Range("K62").Value = 12002
$ws.Range("L62").Value = 550001.5
$ws.Range("M62").Value = -11378
$ws.Range("N62").Value = -551249.5
$ws.Range("H65").Value = 281001.75
$ws.Range("I65").Value = 12002
$ws.Range("J65").Value = 550001.5
$ws.Range("K65").Value = 60010
$ws.Range("L65").Value = 2750007.5
$ws.Range("M65").Value = -56890
$ws.Range("N65").Value = -2756247.5
